$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42 (rows 42..52 shift down to 43..53).
$ws.Rows(42).Insert()

# Populate the newly inserted row 42 with the new data record.
$ws.Range("A42").Value = 10
$ws.Range("B42").Value = "Vega Modelo de Temuco"
$ws.Range("C42").Value = "La Araucanía"
$ws.Range("D42").Value = 44508
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100101
$ws.Range("H42").Value = "Berries"
$ws.Range("I42").Value = 100101001
$ws.Range("J42").Value = "Arándano (blue)"
$ws.Range("K42").Value = "Sin especificar"
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 120
$ws.Range("N42").Value = 4000
$ws.Range("O42").Value = 4000
$ws.Range("P42").Value = 4000
$ws.Range("Q42").Value = "$/kilo"
$ws.Range("R42").Value = "Provincia de Limarí"
$ws.Range("S42").Value = 4000
$ws.Range("T42").Value = 1

# Make sure the date cell keeps the workbook's date number format.
$ws.Range("D42").NumberFormat = $ws.Range("D43").NumberFormat
